$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new sample was recorded for 2026/02/28 (Saturday) at hour 4, inserted
# right before the existing 2026/12/29 block. Insert a whole row at 900 so
# every subsequent row (900-941) shifts down to (901-942), then populate
# the newly inserted row with the new data point.
$ws.Rows("900:900").Insert()

# Force column A to stay plain text (it holds literal "yyyy/mm/dd" strings,
# not real dates, throughout the sheet) - a leading apostrophe tells Excel
# to store the literal text instead of auto-converting it to a date serial.
$ws.Cells.Item(900, 1).Value = "'2026/02/28"
$ws.Cells.Item(900, 1).Style = "Normal"
$ws.Cells.Item(900, 2).Value = "土"
$ws.Cells.Item(900, 3).Value = 4
$ws.Cells.Item(900, 4).Value = 201
